# Added Waitaki flow data
# Append a new "New Zealand (community)" dataset row (Waitaki Flow Dataset)
# to the Community table on Sheet2, and move the selection to I15 (matching
# the author's final cursor position when the edit was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "New Zealand (community)"
$ws.Range("B23").Value = "Waitaki Flow Dataset"
$ws.Range("C23").Value = "Waitaki_flow.xlsx"

[void]$ws.Range("I15").Select()
